$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update Ammar ASLAN's paper title + email ---
$ws.Range("C2").Value = "ammaraslan@test.com"
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:ammaraslan@test.com") | Out-Null
$ws.Range("C2").Style = "Köprü"

$ws.Range("F2").Value = "A HIGH-PERFORMANCE INTEGER LINEAR PROGRAMMING BASED COMPUTATION FOR TRAFFIC SCHEDULES IN IEEE 802.1 TIME SENSITIVE NETWORKS"

# --- Row 3: new participant Selahattin Barış Çelebi ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Selahattin Barış Çelebi"

$ws.Range("C3").Value = "sbariscelebi@test.com"
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:sbariscelebi@test.com") | Out-Null
$ws.Range("C3").Style = "Köprü"

$ws.Range("F3").Value = "LORAWAN KABLOSUZ HABERLEŞME PROTOKOLÜNÜN GÜVENLİK ANALİZİ VE LORAWAN TABANLI IOT CİHAZLARINA KARŞI YAPILAN SALDIRILAR"

# Give the small grey "paper title" font treatment to F2/F3
$titleRange = $ws.Range("F2:F3")
$titleRange.Font.Name = "Open Sans"
$titleRange.Font.Size = 8
$titleRange.Font.Color = 5592405

# Move the active selection like the author left it
$ws.Range("E5").Select() | Out-Null

Write-Host "edit applied"
